$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.055348271424477
$ws.Range("D2").Value = 1.059223765969924
$ws.Range("E2").Value = 1.05172781776135
$ws.Range("F2").Value = 1.068201921739796
$ws.Range("I2").Value = 1.050872577125768
$ws.Range("J2").Value = 1.060355602432158
$ws.Range("K2").Value = 1.061953854772176
$ws.Range("L2").Value = 1.054478516877626
$ws.Range("M2").Value = 1.070907736768231
$ws.Range("N2").Value = 1.0618614290121
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.056593545473205
$ws.Range("D3").Value = 1.060232641963411
$ws.Range("E3").Value = 1.052803067570205
$ws.Range("F3").Value = 1.069443555495317
$ws.Range("I3").Value = 1.051312973854122
$ws.Range("J3").Value = 1.061251089807425
$ws.Range("K3").Value = 1.062776425360184
$ws.Range("L3").Value = 1.055365795996082
$ws.Range("M3").Value = 1.071964245104395
$ws.Range("N3").Value = 1.062758188082153
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.057398757179376
$ws.Range("D4").Value = 1.060884916209817
$ws.Range("E4").Value = 1.053498457051168
$ws.Range("F4").Value = 1.070246899659429
$ws.Range("I4").Value = 1.051596339751782
$ws.Range("J4").Value = 1.0618294440183
$ws.Range("K4").Value = 1.063307536670919
$ws.Range("L4").Value = 1.055938960743286
$ws.Range("M4").Value = 1.07264723478061
$ws.Range("N4").Value = 1.063337363622346
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.057737135940284
$ws.Range("D5").Value = 1.061159005556815
$ws.Range("E5").Value = 1.05379071217558
$ws.Range("F5").Value = 1.070584609505153
$ws.Range("I5").Value = 1.051715084673273
$ws.Range("J5").Value = 1.062072325900816
$ws.Range("K5").Value = 1.063530542856669
$ws.Range("L5").Value = 1.056179690023958
$ws.Range("M5").Value = 1.072934211881747
$ws.Range("N5").Value = 1.063580590424994
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.05779394354847
$ws.Range("D6").Value = 1.061205018977306
$ws.Range("E6").Value = 1.053839778068329
$ws.Range("F6").Value = 1.070641311595459
$ws.Range("I6").Value = 1.051735000090412
$ws.Range("J6").Value = 1.062113091746206
$ws.Range("K6").Value = 1.063567970612191
$ws.Range("L6").Value = 1.056220096138067
$ws.Range("M6").Value = 1.072982387755889
$ws.Range("N6").Value = 1.06362141416256
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.057403279129119
$ws.Range("D7").Value = 1.060888579101159
$ws.Range("E7").Value = 1.05350236251661
$ws.Range("F7").Value = 1.070251412215286
$ws.Range("I7").Value = 1.051597927927394
$ws.Range("J7").Value = 1.061832690430366
$ws.Range("K7").Value = 1.063310517561297
$ws.Range("L7").Value = 1.055942178278183
$ws.Range("M7").Value = 1.072651069976615
$ws.Range("N7").Value = 1.063340614644689
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.055769236030632
$ws.Range("D8").Value = 1.059564832260393
$ws.Range("E8").Value = 1.052091280928493
$ws.Range("F8").Value = 1.068621554680072
$ws.Range("I8").Value = 1.051021743202127
$ws.Range("J8").Value = 1.060658462509244
$ws.Range("K8").Value = 1.062232084754769
$ws.Range("L8").Value = 1.054778577640096
$ws.Range("M8").Value = 1.071264922756942
$ws.Range("N8").Value = 1.062164719185233
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.052885401255118
$ws.Range("D9").Value = 1.057228045462635
$ws.Range("E9").Value = 1.04960187025677
$ws.Range("F9").Value = 1.065748860385298
$ws.Range("I9").Value = 1.049994131317537
$ws.Range("J9").Value = 1.058580926124339
$ws.Range("K9").Value = 1.060322894892113
$ws.Range("L9").Value = 1.052720706197654
$ws.Range("M9").Value = 1.068817352674201
$ws.Range("N9").Value = 1.060084232460365
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.050959683455156
$ws.Range("D10").Value = 1.055667278729996
$ws.Range("E10").Value = 1.047940198580463
$ws.Range("F10").Value = 1.063833139315457
$ws.Range("I10").Value = 1.04930072226064
$ws.Range("J10").Value = 1.057190143455184
$ws.Range("K10").Value = 1.059044049385306
$ws.Range("L10").Value = 1.051343671735485
$ws.Range("M10").Value = 1.06718215448069
$ws.Range("N10").Value = 1.058691474720295
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.050125037052861
$ws.Range("D11").Value = 1.054990736249084
$ws.Range("E11").Value = 1.047220161130898
$ws.Range("F11").Value = 1.063003439549209
$ws.Range("I11").Value = 1.048998476957947
$ws.Range("J11").Value = 1.056586527617759
$ws.Range("K11").Value = 1.058488836825876
$ws.Range("L11").Value = 1.050746162055077
$ws.Range("M11").Value = 1.06647324120019
$ws.Range("N11").Value = 1.058087001679144
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.049814888845443
$ws.Range("D12").Value = 1.054739327943552
$ws.Range("E12").Value = 1.04695262634189
$ws.Range("F12").Value = 1.062695222311306
$ws.Range("I12").Value = 1.048885908591031
$ws.Range("J12").Value = 1.056362105509532
$ws.Range("K12").Value = 1.058282384152149
$ws.Range("L12").Value = 1.050524031189897
$ws.Range("M12").Value = 1.066209787291506
$ws.Range("N12").Value = 1.057862260865782
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.049881422327834
$ws.Range("D13").Value = 1.054793260903959
$ws.Range("E13").Value = 1.047010017156609
$ws.Range("F13").Value = 1.062761337348381
$ws.Range("I13").Value = 1.04891006853044
$ws.Range("J13").Value = 1.056410254458265
$ws.Range("K13").Value = 1.058326679030933
$ws.Range("L13").Value = 1.050571687609561
$ws.Range("M13").Value = 1.066266305060507
$ws.Range("N13").Value = 1.057910478191545
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.050099402642282
$ws.Range("D14").Value = 1.054969957022521
$ws.Range("E14").Value = 1.047198048283821
$ws.Range("F14").Value = 1.062977962840216
$ws.Range("I14").Value = 1.048989178169761
$ws.Range("J14").Value = 1.05656798116075
$ws.Range("K14").Value = 1.058471775924286
$ws.Range("L14").Value = 1.050727804519957
$ws.Range("M14").Value = 1.066451466724158
$ws.Range("N14").Value = 1.058068428884039
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.050233691055616
$ws.Range("D15").Value = 1.05507881065325
$ws.Range("E15").Value = 1.047313889688691
$ws.Range("F15").Value = 1.063111428923022
$ws.Range("I15").Value = 1.049037880303471
$ws.Range("J15").Value = 1.056665133585325
$ws.Range("K15").Value = 1.058561145429181
$ws.Range("L15").Value = 1.050823968169685
$ws.Range("M15").Value = 1.066565533376527
$ws.Range("N15").Value = 1.058165719276202
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.051015059079194
$ws.Range("D16").Value = 1.055712163273342
$ws.Range("E16").Value = 1.047987973887013
$ws.Range("F16").Value = 1.063888199729519
$ws.Range("I16").Value = 1.049320739147884
$ws.Range("J16").Value = 1.057230173803924
$ws.Range("K16").Value = 1.059080865993839
$ws.Range("L16").Value = 1.051383300076452
$ws.Range("M16").Value = 1.067229184353133
$ws.Range("N16").Value = 1.058731561916722
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.051504974513583
$ws.Range("D17").Value = 1.056109254340857
$ws.Range("E17").Value = 1.048410667547244
$ws.Range("F17").Value = 1.064375397413792
$ws.Range("I17").Value = 1.04949763407254
$ws.Range("J17").Value = 1.057584232753717
$ws.Range("K17").Value = 1.059406479564711
$ws.Range("L17").Value = 1.051733819441685
$ws.Range("M17").Value = 1.067645242623809
$ws.Range("N17").Value = 1.059086123670835
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.05179065715108
$ws.Range("D18").Value = 1.056340801498521
$ws.Range("E18").Value = 1.04865716733666
$ws.Range("F18").Value = 1.064659554411514
$ws.Range("I18").Value = 1.049600621529895
$ws.Range("J18").Value = 1.057790614642505
$ws.Range("K18").Value = 1.059596263335551
$ws.Range("L18").Value = 1.051938151512132
$ws.Range("M18").Value = 1.067887839405164
$ws.Range("N18").Value = 1.059292798645577
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.051888054635397
$ws.Range("D19").Value = 1.056419741396784
$ws.Range("E19").Value = 1.048741208936827
$ws.Range("F19").Value = 1.064756441775641
$ws.Range("I19").Value = 1.049635705004085
$ws.Range("J19").Value = 1.057860962754088
$ws.Range("K19").Value = 1.05966095082295
$ws.Range("L19").Value = 1.052007803188088
$ws.Range("M19").Value = 1.067970544640329
$ws.Range("N19").Value = 1.059363246659548
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.051452419188869
$ws.Range("D20").Value = 1.056066657416142
$ws.Range("E20").Value = 1.048365321735318
$ws.Range("F20").Value = 1.064323127510193
$ws.Range("I20").Value = 1.049478674828146
$ws.Range("J20").Value = 1.057546259530221
$ws.Range("K20").Value = 1.059371558922231
$ws.Range("L20").Value = 1.051696224454199
$ws.Range("M20").Value = 1.067600612109902
$ws.Range("N20").Value = 1.059048096521006
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.050035216306477
$ws.Range("D21").Value = 1.054917927481551
$ws.Range("E21").Value = 1.047142680052955
$ws.Range("F21").Value = 1.062914172877137
$ws.Range("I21").Value = 1.04896589067268
$ws.Range("J21").Value = 1.056521540457511
$ws.Range("K21").Value = 1.058429054650706
$ws.Range("L21").Value = 1.050681837220304
$ws.Range("M21").Value = 1.066396944891849
$ws.Range("N21").Value = 1.058021922229674
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.049143449079649
$ws.Range("D22").Value = 1.054195036138423
$ws.Range("E22").Value = 1.046373487953305
$ws.Range("F22").Value = 1.062028133370018
$ws.Range("I22").Value = 1.048641740740234
$ws.Range("J22").Value = 1.055876029393338
$ws.Range("K22").Value = 1.05783517931189
$ws.Range("L22").Value = 1.050042956171885
$ws.Range("M22").Value = 1.065639386863071
$ws.Range("N22").Value = 1.05737549446575
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.049616260754123
$ws.Range("D23").Value = 1.054578315815282
$ws.Range("E23").Value = 1.046781296440919
$ws.Range("F23").Value = 1.062497857081511
$ws.Range("I23").Value = 1.048813744339002
$ws.Range("J23").Value = 1.056218344339964
$ws.Range("K23").Value = 1.058150126404571
$ws.Range("L23").Value = 1.050381743673751
$ws.Range("M23").Value = 1.066041056025948
$ws.Range("N23").Value = 1.057718295538864
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.051476166914931
$ws.Range("D24").Value = 1.056085905347131
$ws.Range("E24").Value = 1.048385811711174
$ws.Range("F24").Value = 1.064346746082731
$ws.Range("I24").Value = 1.049487242288937
$ws.Range("J24").Value = 1.05756341841157
$ws.Range("K24").Value = 1.059387338493434
$ws.Range("L24").Value = 1.051713212381379
$ws.Range("M24").Value = 1.067620778974378
$ws.Range("N24").Value = 1.059065279769934
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.053631486032742
$ws.Range("D25").Value = 1.057832666112645
$ws.Range("E25").Value = 1.050245798268588
$ws.Range("F25").Value = 1.066491614960955
$ws.Range("I25").Value = 1.050261257915347
$ws.Range("J25").Value = 1.059119025736851
$ws.Range("K25").Value = 1.060817525656903
$ws.Range("L25").Value = 1.05325360999982
$ws.Range("M25").Value = 1.069450713458159
$ws.Range("N25").Value = 1.060623096236047
